$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A8:XFD8").Select()
$win = $excel.ActiveWindow
$win.Zoom = 120
$win.ScrollRow = 4
$win.ScrollColumn = 1
Write-Host "done"
$ws2 = $wb.Worksheets.Item(2)
$ws2.Activate()
$ws2.Range("F4").Select()
